$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $origStyle = $r.Style
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = $origStyle
}

Set-TextValue "D2" "70.598.46"
$ws.Range("E2").Value = "  -0.22%  "
Set-TextValue "D3" "3.647.66"
$ws.Range("E3").Value = "  +5.11%  "
$ws.Range("E4").Value = "  +0.34%  "
Set-TextValue "D5" "577.38"
$ws.Range("E5").Value = "  -1.79%  "
Set-TextValue "D6" "176.23"
$ws.Range("E6").Value = "  -2.00%  "
Set-TextValue "D7" "3.640.47"
$ws.Range("E7").Value = "  +5.41%  "
Set-TextValue "D8" "0.613"
$ws.Range("E8").Value = "  +1.75%  "
$ws.Range("E9").Value = "  +0.20%  "
$ws.Range("E10").Value = "  -4.44%  "
Set-TextValue "D11" "6.84"
$ws.Range("E11").Value = "  +23.82%  "
Set-TextValue "D12" "0.605"
$ws.Range("E12").Value = "  +1.94%  "
Set-TextValue "D13" "48.63"
$ws.Range("E13").Value = "  -1.64%  "
Set-TextValue "D14" "0.0000288"
$ws.Range("E14").Value = "  -0.15%  "
Set-TextValue "D15" "4.234.53"
$ws.Range("E15").Value = "  +5.89%  "
Set-TextValue "D16" "671.52"
$ws.Range("E16").Value = "  -3.45%  "
Set-TextValue "D17" "8.89"
$ws.Range("E17").Value = "  +1.39%  "
Set-TextValue "D18" "3.642.23"
$ws.Range("E18").Value = "  +5.89%  "
Set-TextValue "D19" "70.724.15"
$ws.Range("E19").Value = "  +0.33%  "
$ws.Range("E20").Value = "  +0.46%  "
Set-TextValue "D21" "17.79"
$ws.Range("E21").Value = "  -0.56%  "
Set-TextValue "D22" "11.42"
$ws.Range("E22").Value = "  -1.06%  "
Set-TextValue "D23" "0.930"
$ws.Range("E23").Value = "  +2.04%  "
Set-TextValue "D24" "17.11"
$ws.Range("E24").Value = "  -0.58%  "
Set-TextValue "D25" "100.53"
$ws.Range("E25").Value = "  -1.31%  "
$ws.Range("E26").Value = "  -1.58%  "
Set-TextValue "D27" "2.79"
$ws.Range("E27").Value = "  +2.89%  "
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("E29").Value = "  +2.39%  "
Set-TextValue "D30" "34.97"
$ws.Range("E30").Value = "  +2.91%  "
$ws.Range("E31").Value = "  -0.44%  "
$ws.Range("E32").Value = "  +1.14%  "
$ws.Range("E33").Value = "  -4.67%  "
Set-TextValue "D34" "7.31"
$ws.Range("E34").Value = "  +0.83%  "
$ws.Range("E35").Value = "  +1.20%  "
Set-TextValue "D36" "585.15"
$ws.Range("E36").Value = "  +1.39%  "
Set-TextValue "D37" "11.05"
$ws.Range("E37").Value = "  -0.67%  "
$ws.Range("E38").Value = "  +2.63%  "
Set-TextValue "D39" "58.15"
$ws.Range("E39").Value = "  -1.33%  "
$ws.Range("E40").Value = "  +0.07%  "
Set-TextValue "D41" "3.581.16"
$ws.Range("E41").Value = "  -0.66%  "
Set-TextValue "D42" "0.0454"
$ws.Range("E42").Value = "  +6.57%  "
Set-TextValue "D43" "0.142"
$ws.Range("E43").Value = "  +0.73%  "
Set-TextValue "D44" "0.344"
$ws.Range("E44").Value = "  +1.22%  "
Set-TextValue "D45" "34.80"
$ws.Range("E45").Value = "  -2.41%  "
$ws.Range("E46").Value = "  -1.01%  "
$ws.Range("E47").Value = "  -1.47%  "
Set-TextValue "D48" "2.94"
$ws.Range("E48").Value = "  +8.60%  "
Set-TextValue "D49" "0.133"
$ws.Range("E49").Value = "  +1.98%  "
Set-TextValue "D50" "135.62"
$ws.Range("E50").Value = "  +1.60%  "
$ws.Range("E51").Value = "  +6.43%  "
